# "Aggiornamento audio + logo owof new style + ore"
# Adds new time-tracking entries (rows 46-48) to the "Prot. 2.0" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prot. 2.0")
$ws.Activate()

# Row 46: 09/03/2025 (serial 45725) - Programmazione
$ws.Range("A46").Value = 45725
$ws.Range("B46").Value = "Programmazione"
$ws.Range("C46").Value = "Risoluzione con Mattia di una serie di problemi vari"
$ws.Range("D46").Value = 2/24

# Row 47: Gestione
$ws.Range("B47").Value = "Gestione"
$ws.Range("C47").Value = "Testing e risoluzione piccoli problemi"
$ws.Range("D47").Value = 0.5/24

# Row 48: Grafica
$ws.Range("B48").Value = "Grafica"
$ws.Range("C48").Value = "Gestione e impostazione scena d'avvio gioco"
$ws.Range("D48").Value = 0.5/24

# Note detail added last for row 46
$ws.Range("E46").Value = "Risoluzione problema slot bianchi di testo, settaggio scena avvio, settaggio opzione reset gioco, overlay tasti, risoluzione problema biblioteca"
$ws.Range("E46").WrapText = $true
$ws.Rows.Item(46).RowHeight = 43.2

$ws.Range("E37").Select()
